# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (fund holdings detail) positioned
#    between "2021-Q4" and "总计".
# 2) Insert a new summary row at the top of "总计" for the 2022-Q1 quarter
#    and renumber the index column (A) for the existing rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create & position the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add()
$new.Name = "2022-Q1"
$new.Move($wb.Worksheets.Item("总计"))

# Re-fetch by name: some host object references go stale after Move().
$ws = $wb.Worksheets.Item("2022-Q1")

$totals = $wb.Worksheets.Item("总计")

# --- header row -------------------------------------------------------
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# --- index column (A2:A14) --------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12

# --- fund code column (B2:B14) -- keep as text (leading zeros) --------
$codes = $ws.Range("B2:B14")
$codes.NumberFormat = "@"
$ws.Range("B2").Value = "002560"
$ws.Range("B3").Value = "005777"
$ws.Range("B4").Value = "000654"
$ws.Range("B5").Value = "004423"
$ws.Range("B6").Value = "320022"
$ws.Range("B7").Value = "012491"
$ws.Range("B8").Value = "008328"
$ws.Range("B9").Value = "008961"
$ws.Range("B10").Value = "012492"
$ws.Range("B11").Value = "005901"
$ws.Range("B12").Value = "001351"
$ws.Range("B13").Value = "005902"
$ws.Range("B14").Value = "010355"
$codes.Style = "Normal"

# --- fund name column (C2:C14) -----------------------------------------
$ws.Range("C2").Value = "诺安和鑫灵活配置混合"
$ws.Range("C3").Value = "广发科技动力股票"
$ws.Range("C4").Value = "华商新锐产业灵活配置混合"
$ws.Range("C5").Value = "华商研究精选灵活配置混合"
$ws.Range("C6").Value = "诺安研究精选股票"
$ws.Range("C7").Value = "华商核心引力混合型证券投资基金A"
$ws.Range("C8").Value = "诺安新兴产业混合"
$ws.Range("C9").Value = "华商科技创新混合"
$ws.Range("C10").Value = "华商核心引力混合型证券投资基金C"
$ws.Range("C11").Value = "诺安汇利灵活配置混合A"
$ws.Range("C12").Value = "诺安中证500指数增强A"
$ws.Range("C13").Value = "诺安汇利灵活配置混合C"
$ws.Range("C14").Value = "诺安中证500指数增强C"

# --- D:G numeric-looking columns -- stored as text, like the sibling
#     quarter sheets (2021-Q2 / 2021-Q4) ---------------------------------
$numericText = $ws.Range("D2:G14")
$numericText.NumberFormat = "@"

$ws.Range("D2").Value = "33.85"
$ws.Range("E2").Value = "93.22"
$ws.Range("F2").Value = "7.10"
$ws.Range("G2").Value = "2.4034"

$ws.Range("D3").Value = "17.47"
$ws.Range("E3").Value = "92.65"
$ws.Range("F3").Value = "5.32"
$ws.Range("G3").Value = "0.9294"

$ws.Range("D4").Value = "13.02"
$ws.Range("E4").Value = "81.79"
$ws.Range("F4").Value = "2.51"
$ws.Range("G4").Value = "0.3268"

$ws.Range("D5").Value = "9.97"
$ws.Range("E5").Value = "82.17"
$ws.Range("F5").Value = "2.51"
$ws.Range("G5").Value = "0.2502"

$ws.Range("D6").Value = "8.35"
$ws.Range("E6").Value = "93.85"
$ws.Range("F6").Value = "2.72"
$ws.Range("G6").Value = "0.2271"

$ws.Range("D7").Value = "5.05"
$ws.Range("E7").Value = "83.13"
$ws.Range("F7").Value = "2.51"
$ws.Range("G7").Value = "0.1268"

$ws.Range("D8").Value = "4.83"
$ws.Range("E8").Value = "68.49"
$ws.Range("F8").Value = "2.36"
$ws.Range("G8").Value = "0.1140"

$ws.Range("D9").Value = "2.87"
$ws.Range("E9").Value = "88.86"
$ws.Range("F9").Value = "2.51"
$ws.Range("G9").Value = "0.0720"

$ws.Range("D10").Value = "0.51"
$ws.Range("E10").Value = "83.13"
$ws.Range("F10").Value = "2.51"
$ws.Range("G10").Value = "0.0128"

$ws.Range("D11").Value = "0.08"
$ws.Range("E11").Value = "86.88"
$ws.Range("F11").Value = "4.28"
$ws.Range("G11").Value = "0.0034"

$ws.Range("D12").Value = "0.48"
$ws.Range("E12").Value = "94.42"
$ws.Range("F12").Value = "0.65"
$ws.Range("G12").Value = "0.0031"

$ws.Range("D13").Value = "0.02"
$ws.Range("E13").Value = "86.88"
$ws.Range("F13").Value = "4.28"
$ws.Range("G13").Value = "0.0009"

$ws.Range("D14").Value = "0.04"
$ws.Range("E14").Value = "94.42"
$ws.Range("F14").Value = "0.65"
$ws.Range("G14").Value = "0.0003"

$numericText.Style = "Normal"

# --- rank column (H2:H14) -- numeric ------------------------------------
$ws.Range("H2").Value = 7
$ws.Range("H3").Value = 7
$ws.Range("H4").Value = 8
$ws.Range("H5").Value = 8
$ws.Range("H6").Value = 6
$ws.Range("H7").Value = 8
$ws.Range("H8").Value = 10
$ws.Range("H9").Value = 9
$ws.Range("H10").Value = 8
$ws.Range("H11").Value = 10
$ws.Range("H12").Value = 3
$ws.Range("H13").Value = 10
$ws.Range("H14").Value = 3

# --- formatting: reuse the existing "header" / "index column" styles
#     already present on the "总计" sheet instead of re-building fonts /
#     borders from scratch (keeps the style table minimal, matches the
#     sibling sheets). ----------------------------------------------------
$totals.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$totals.Range("A2").Copy()
$ws.Range("A2:A14").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 summary row on "总计"
# ---------------------------------------------------------------------
$totals.Rows(2).Insert()
$totals.Range("A2:D2").Style = "Normal"

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 13
$totals.Range("D2").Value = 4.47

# renumber the index column for the rows that shifted down
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2

$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

$wb.Worksheets.Item("2021-Q2").Activate()
